$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the text labels ("1号".."5号") in column A with plain numeric values
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5

# Update visit-count values in column B
$ws.Range("B4").Value = 100
$ws.Range("B6").Value = 50

# Move the active selection from B6 to C6
$ws.Range("C6").Select()
